# Update citation DOI cells AP10:AP38 of the compilation sheet so that each
# one references the cell directly above it, instead of holding an
# independent (and incorrectly incrementing) DOI string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AP10 references AP9 directly.
$ws.Range("AP10").Formula = "=AP9"

# AP11:AP38 each reference the cell directly above them (AP11=AP10,
# AP12=AP11, ... AP38=AP37).
$ws.Range("AP11:AP38").FormulaR1C1 = "=R[-1]C"

# The old AP10:AP38 hyperlink (which pointed at a set of now-removed,
# individually-numbered DOI strings) is no longer meaningful now that those
# cells are formulas mirroring AP9, so drop it. Only the original AP9
# hyperlink should remain.
$keepAddress = '$AP$9'
for ($i = $ws.Hyperlinks.Count; $i -ge 1; $i--) {
    $hl = $ws.Hyperlinks.Item($i)
    $addr = $hl.Range.Address()
    if ($addr -ne $keepAddress) {
        $hl.Delete() | Out-Null
    }
}

# Update the selection shown when the file is opened to match AP10:AP38,
# with AP10 as the active cell.
$ws.Range("AP10:AP38").Select() | Out-Null
